# Update the date values on Sheet2 (add one month's worth / +100 to the
# yyyymmdd-as-integer values stored in column B).
$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Sheet2")

$rows = @(2, 5, 8, 11, 14, 17, 20, 23, 26, 29)
foreach ($r in $rows) {
    $cell = $ws2.Range("B$r")
    $cell.Value = $cell.Value2 + 100
}

# Move the selection/active cell on Sheet2 back to the top and select H12.
$ws2.Activate() | Out-Null
$ws2.Range("A1").Select() | Out-Null
$ws2.Range("H12").Select() | Out-Null

# Move the selection/active cell on Sheet1 to A9.
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Activate() | Out-Null
$ws1.Range("A9").Select() | Out-Null

# Restore Sheet2 as the active sheet (it was tabSelected in the workbook).
$ws2.Activate() | Out-Null
